# Weekly fruit/vegetable data update ("Fruta / hortaliza, semanal").
# Two new daily price records are inserted at the top of the data block
# (rows 804-805), pushing the existing records (old rows 804-901) down by
# two rows so they become rows 806-903. The sheet's used range grows from
# A1:R901 to A1:R903.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 804; everything from the old row 804 onward
# shifts down to make room (old 804 -> 806, old 805 -> 807, ... old 901 -> 903).
$ws.Rows.Item(804).Resize(2).Insert()

# ---- New row 804: Lechuga Conconina(o), Región del Maule ----
$ws.Cells.Item(804, 1).Value = 7
$ws.Cells.Item(804, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(804, 3).Value = "Ñuble"
$ws.Cells.Item(804, 4).Value = 44918
$ws.Cells.Item(804, 5).Value = 16
$ws.Cells.Item(804, 6).Value = 100112033
$ws.Cells.Item(804, 7).Value = "Lechuga"
$ws.Cells.Item(804, 8).Value = "Conconina(o)"
$ws.Cells.Item(804, 9).Value = "Primera"
$ws.Cells.Item(804, 10).Value = 120
$ws.Cells.Item(804, 11).Value = 4500
$ws.Cells.Item(804, 12).Value = 5000
$ws.Cells.Item(804, 13).Value = 4750
$ws.Cells.Item(804, 14).Value = "`$/caja 10 unidades"
$ws.Cells.Item(804, 15).Value = "Región del Maule"
$ws.Cells.Item(804, 16).Value = 475
$ws.Cells.Item(804, 17).Value = 10
$ws.Cells.Item(804, 18).Value = "Hortaliza"

# ---- New row 805: Lechuga Escarola, Región del Maule ----
$ws.Cells.Item(805, 1).Value = 7
$ws.Cells.Item(805, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(805, 3).Value = "Ñuble"
$ws.Cells.Item(805, 4).Value = 44918
$ws.Cells.Item(805, 5).Value = 16
$ws.Cells.Item(805, 6).Value = 100112033
$ws.Cells.Item(805, 7).Value = "Lechuga"
$ws.Cells.Item(805, 8).Value = "Escarola"
$ws.Cells.Item(805, 9).Value = "Primera"
$ws.Cells.Item(805, 10).Value = 120
$ws.Cells.Item(805, 11).Value = 6000
$ws.Cells.Item(805, 12).Value = 6500
$ws.Cells.Item(805, 13).Value = 6250
$ws.Cells.Item(805, 14).Value = "`$/caja 15 unidades"
$ws.Cells.Item(805, 15).Value = "Región del Maule"
$ws.Cells.Item(805, 16).Value = 417
$ws.Cells.Item(805, 17).Value = 15
$ws.Cells.Item(805, 18).Value = "Hortaliza"
